# Apply Week 17 stat updates to the Panthers "Players Data" workbook.
# Sheet "Rushing" holds rushing attempt splits; sheet "Receiving" holds
# target/completion splits. Only raw counting stats are bumped up to
# reflect the newly logged week's numbers - no structural changes.

$wb = $excel.ActiveWorkbook

$wsRushing = $wb.Worksheets.Item("Rushing")
$wsReceiving = $wb.Worksheets.Item("Receiving")

# --- Rushing sheet updates ---
# Row 2: C.Newton
$wsRushing.Range("C2").Value = 25

# Row 4: S.Darnold
$wsRushing.Range("C4").Value = 2
$wsRushing.Range("E4").Value = 1

# Row 5: C.Hubbard
$wsRushing.Range("C5").Value = 43
$wsRushing.Range("D5").Value = 24
$wsRushing.Range("E5").Value = 9
$wsRushing.Range("F5").Value = 13

# Row 7: A.Abdullah
$wsRushing.Range("C7").Value = 17
$wsRushing.Range("D7").Value = 16

# --- Receiving sheet updates ---
# Row 2: C.Hubbard
$wsReceiving.Range("C2").Value = 21
$wsReceiving.Range("D2").Value = 13

# Row 4: A.Abdullah
$wsReceiving.Range("C4").Value = 35
$wsReceiving.Range("D4").Value = 25
$wsReceiving.Range("G4").Value = 5
$wsReceiving.Range("H4").Value = 3

# Row 6: R.Anderson
$wsReceiving.Range("C6").Value = 79
$wsReceiving.Range("D6").Value = 42

# Row 7: Dj.Moore
$wsReceiving.Range("C7").Value = 114
$wsReceiving.Range("D7").Value = 71
$wsReceiving.Range("E7").Value = 39
$wsReceiving.Range("F7").Value = 15

# Row 8: T.Marshall
$wsReceiving.Range("C8").Value = 24
$wsReceiving.Range("D8").Value = 16

# Row 10: S.Smith
$wsReceiving.Range("C10").Value = 9

# Row 15: I.Thomas
$wsReceiving.Range("C15").Value = 26
$wsReceiving.Range("D15").Value = 15
$wsReceiving.Range("E15").Value = 4
$wsReceiving.Range("F15").Value = 3
